$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update praclen column (I) values from 4 to 5 for existing rows 2-5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Add new row 6 with trial data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 61
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_1"

# Update selection to K6
$ws.Range("K6").Select()
